# Auto-generated-style script: applies market-data refresh values
# to the profit calculation columns (H-N) across all 8 sheets.
$wb = $excel.ActiveWorkbook

# ALC!row6
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(6, 8).Value = 2041548.5
$ws.Cells.Item(6, 9).Value = 4762280
$ws.Cells.Item(6, 11).Value = 14286840
$ws.Cells.Item(6, 13).Value = -14286728

# ALC!row8
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(8, 8).Value = 8547268
$ws.Cells.Item(8, 9).Value = 9615639
$ws.Cells.Item(8, 10).Value = 300
$ws.Cells.Item(8, 11).Value = 28846917
$ws.Cells.Item(8, 12).Value = 900
$ws.Cells.Item(8, 13).Value = -28846778
$ws.Cells.Item(8, 14).Value = -1178

# ALC!row19
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(19, 8).Value = 565.3333
$ws.Cells.Item(19, 9).Value = 253.57143
$ws.Cells.Item(19, 10).Value = 693.7059
$ws.Cells.Item(19, 11).Value = 253.57143
$ws.Cells.Item(19, 12).Value = 693.7059
$ws.Cells.Item(19, 13).Value = -78.57142999999999
$ws.Cells.Item(19, 14).Value = -1043.7059

# ALC!row96
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(96, 8).Value = 2564417.5
$ws.Cells.Item(96, 9).Value = 3663235.2
$ws.Cells.Item(96, 10).Value = 509.66666
$ws.Cells.Item(96, 11).Value = 10989705.6
$ws.Cells.Item(96, 12).Value = 1528.99998
$ws.Cells.Item(96, 13).Value = -10988332.6
$ws.Cells.Item(96, 14).Value = -4274.999980000001

# ALC!row129
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(129, 8).Value = 1204.1831
$ws.Cells.Item(129, 9).Value = 868.1429000000001
$ws.Cells.Item(129, 10).Value = 1240.9375
$ws.Cells.Item(129, 11).Value = 2604.4287
$ws.Cells.Item(129, 12).Value = 3722.8125
$ws.Cells.Item(129, 13).Value = 2395.5713
$ws.Cells.Item(129, 14).Value = -13722.8125

# ALC!row137
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(137, 8).Value = 2070.923
$ws.Cells.Item(137, 9).Value = 5109
$ws.Cells.Item(137, 10).Value = 1159.5
$ws.Cells.Item(137, 11).Value = 15327
$ws.Cells.Item(137, 12).Value = 3478.5
$ws.Cells.Item(137, 13).Value = -12777
$ws.Cells.Item(137, 14).Value = -8578.5

# ARM!row32
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 7929.4604
$ws.Cells.Item(32, 9).Value = 7110.71
$ws.Cells.Item(32, 10).Value = 16000
$ws.Cells.Item(32, 11).Value = 7110.71
$ws.Cells.Item(32, 12).Value = 16000
$ws.Cells.Item(32, 13).Value = -6823.71
$ws.Cells.Item(32, 14).Value = -16574

# ARM!row122
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(122, 8).Value = 1082.7059
$ws.Cells.Item(122, 9).Value = 872.5454999999999
$ws.Cells.Item(122, 10).Value = 1468
$ws.Cells.Item(122, 11).Value = 2617.6365
$ws.Cells.Item(122, 12).Value = 4404
$ws.Cells.Item(122, 13).Value = -167.6364999999996
$ws.Cells.Item(122, 14).Value = -9304

# BSM!row5
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(5, 8).Value = 738.5
$ws.Cells.Item(5, 9).Value = 888
$ws.Cells.Item(5, 10).Value = 469.4
$ws.Cells.Item(5, 11).Value = 888
$ws.Cells.Item(5, 12).Value = 469.4
$ws.Cells.Item(5, 13).Value = -775
$ws.Cells.Item(5, 14).Value = -695.4

# BSM!row80
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(80, 8).Value = 239.34483
$ws.Cells.Item(80, 9).Value = 217.4375
$ws.Cells.Item(80, 10).Value = 266.30768
$ws.Cells.Item(80, 11).Value = 217.4375
$ws.Cells.Item(80, 12).Value = 266.30768
$ws.Cells.Item(80, 13).Value = 780.5625
$ws.Cells.Item(80, 14).Value = -2262.30768

# BSM!row83
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(83, 8).Value = 239.34483
$ws.Cells.Item(83, 9).Value = 217.4375
$ws.Cells.Item(83, 10).Value = 266.30768
$ws.Cells.Item(83, 11).Value = 1087.1875
$ws.Cells.Item(83, 12).Value = 1331.5384
$ws.Cells.Item(83, 13).Value = 3904.8125
$ws.Cells.Item(83, 14).Value = -11315.5384

# BSM!row86
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(86, 8).Value = 1543.3334
$ws.Cells.Item(86, 9).Value = 1500
$ws.Cells.Item(86, 10).Value = 1546.4286
$ws.Cells.Item(86, 11).Value = 1500
$ws.Cells.Item(86, 12).Value = 1546.4286
$ws.Cells.Item(86, 13).Value = -377
$ws.Cells.Item(86, 14).Value = -3792.4286

# BSM!row89
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(89, 8).Value = 1543.3334
$ws.Cells.Item(89, 9).Value = 1500
$ws.Cells.Item(89, 10).Value = 1546.4286
$ws.Cells.Item(89, 11).Value = 7500
$ws.Cells.Item(89, 12).Value = 7732.143
$ws.Cells.Item(89, 13).Value = -1884
$ws.Cells.Item(89, 14).Value = -18964.143

# BSM!row134
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = 5890.846
$ws.Cells.Item(134, 9).Value = 2407.5625
$ws.Cells.Item(134, 10).Value = 8314
$ws.Cells.Item(134, 11).Value = 7222.6875
$ws.Cells.Item(134, 12).Value = 24942
$ws.Cells.Item(134, 13).Value = -4687.6875
$ws.Cells.Item(134, 14).Value = -30012

# CRP!row22
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(22, 8).Value = 724.85
$ws.Cells.Item(22, 9).Value = 583.46155
$ws.Cells.Item(22, 10).Value = 987.4286
$ws.Cells.Item(22, 11).Value = 583.46155
$ws.Cells.Item(22, 12).Value = 987.4286
$ws.Cells.Item(22, 13).Value = -233.46155
$ws.Cells.Item(22, 14).Value = -1687.4286

# CRP!row31
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 5850542.5
$ws.Cells.Item(31, 9).Value = 2011.5333
$ws.Cells.Item(31, 10).Value = 12348911
$ws.Cells.Item(31, 11).Value = 2011.5333
$ws.Cells.Item(31, 12).Value = 12348911
$ws.Cells.Item(31, 13).Value = -1716.5333
$ws.Cells.Item(31, 14).Value = -12349501

# CRP!row34
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(34, 8).Value = 5850542.5
$ws.Cells.Item(34, 9).Value = 2011.5333
$ws.Cells.Item(34, 10).Value = 12348911
$ws.Cells.Item(34, 11).Value = 2011.5333
$ws.Cells.Item(34, 12).Value = 12348911
$ws.Cells.Item(34, 13).Value = -1809.5333
$ws.Cells.Item(34, 14).Value = -12349315

# CRP!row134
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(134, 8).Value = 2076.72
$ws.Cells.Item(134, 9).Value = 1601.6
$ws.Cells.Item(134, 10).Value = 2789.4
$ws.Cells.Item(134, 11).Value = 4804.799999999999
$ws.Cells.Item(134, 12).Value = 8368.200000000001
$ws.Cells.Item(134, 13).Value = -2269.799999999999
$ws.Cells.Item(134, 14).Value = -13438.2

# CUL!row6
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(6, 8).Value = 635.75
$ws.Cells.Item(6, 9).Value = 177.125
$ws.Cells.Item(6, 10).Value = 1553
$ws.Cells.Item(6, 11).Value = 531.375
$ws.Cells.Item(6, 12).Value = 4659
$ws.Cells.Item(6, 13).Value = -418.375
$ws.Cells.Item(6, 14).Value = -4885

# CUL!row105
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(105, 8).Value = 10380.2
$ws.Cells.Item(105, 10).Value = 10975.111
$ws.Cells.Item(105, 12).Value = 32925.333
$ws.Cells.Item(105, 14).Value = -38167.333

# CUL!row113
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(113, 8).Value = 467.39474
$ws.Cells.Item(113, 9).Value = 407.95834
$ws.Cells.Item(113, 11).Value = 1223.87502
$ws.Cells.Item(113, 13).Value = 946.1249800000001

# GSM!row111
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 14).ClearContents()

# GSM!row122
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(122, 8).Value = 3984.9375
$ws.Cells.Item(122, 9).Value = 3233.25
$ws.Cells.Item(122, 10).Value = 6240
$ws.Cells.Item(122, 11).Value = 9699.75
$ws.Cells.Item(122, 12).Value = 18720
$ws.Cells.Item(122, 13).Value = -7249.75
$ws.Cells.Item(122, 14).Value = -23620

# GSM!row132
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(132, 8).Value = 1265033.8
$ws.Cells.Item(132, 9).Value = 2316533.5
$ws.Cells.Item(132, 10).Value = 3234.1333
$ws.Cells.Item(132, 11).Value = 6949600.5
$ws.Cells.Item(132, 12).Value = 9702.3999
$ws.Cells.Item(132, 13).Value = -6947070.5
$ws.Cells.Item(132, 14).Value = -14762.3999

# LTW!row40
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(40, 8).Value = 76926430
$ws.Cells.Item(40, 9).Value = 111113290
$ws.Cells.Item(40, 10).Value = 5997.5
$ws.Cells.Item(40, 11).Value = 111113290
$ws.Cells.Item(40, 12).Value = 5997.5
$ws.Cells.Item(40, 13).Value = -111113154
$ws.Cells.Item(40, 14).Value = -6269.5

# LTW!row136
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(136, 8).Value = 2916.6667
$ws.Cells.Item(136, 9).Value = 2800
$ws.Cells.Item(136, 10).Value = 2940
$ws.Cells.Item(136, 11).Value = 8400
$ws.Cells.Item(136, 12).Value = 8820
$ws.Cells.Item(136, 13).Value = -5850
$ws.Cells.Item(136, 14).Value = -13920

# WVR!row11
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(11, 8).Value = 657001.9
$ws.Cells.Item(11, 9).Value = 1250250
$ws.Cells.Item(11, 10).Value = 63753.75
$ws.Cells.Item(11, 11).Value = 1250250
$ws.Cells.Item(11, 12).Value = 63753.75
$ws.Cells.Item(11, 13).Value = -1250108
$ws.Cells.Item(11, 14).Value = -64037.75

# WVR!row132
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(132, 8).Value = 2164.4614
$ws.Cells.Item(132, 9).Value = 1870.9
$ws.Cells.Item(132, 11).Value = 5612.700000000001
$ws.Cells.Item(132, 13).Value = -3082.700000000001

